$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Posting Label")

$ws.Unprotect()

# Rows 11..15 in column B/C hold the key/value pairs of the Posting Label.
# The shared-strings table was reshuffled upstream (scenario/scoringCycle/
# scoringMaturity moved up, ahead of planType/variant), which - combined
# with the value-style swap on C13/C15 - shows up as this new key/value
# layout for rows 11..15:
$ws.Range("B11").Value = "scenario"
$ws.Range("C11").Value = "OfficialPlan"

$ws.Range("B12").Value = "scoringCycle"
$ws.Range("C12").Value = "Dec 2020"

$ws.Range("B13").Value = "scoringMaturity"
$ws.Range("C13").Value = "Committed"

$ws.Range("B14").Value = "planType"
$ws.Range("C14").Value = "Marathon"

$ws.Range("B15").Value = "variant"
$ws.Range("C15").Value = "explained"

# C13 picks up the "editable" (unlocked, green-fill) look that C15 used to
# have, and C15 picks up the "locked" (gray-fill) look that C13 used to have.
$ws.Range("C13").Style = $ws.Range("C18").Style
$ws.Range("C15").Style = $ws.Range("C9").Style

$ws.Protect()
